$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- workbook view: tabRatio 986 -> 982 -------------------------------------------------
$wb.Windows.Item(1).TabRatio = 982

# --- touch column C on rows 109..135 (adds an empty, default-styled cell) --------------
for ($r = 109; $r -le 135; $r++) {
  $ws.Cells.Item($r, 3).NumberFormat = "General"
}

# --- new font + style used for the appended "Species" rows (col B) ---------------------
# (font: Arial 10, family "swiss" i.e. family=2 ; numFmt General ; wrap text)
$speciesStyle = $ws.Cells.Item(218, 2).Style

# --- append the new "Species | Tierart" concept block -----------------------------------
$species = @(
  "cattle|Rind",
  "pig|Schwein",
  "sheep|Schaf",
  "goat|Ziege",
  "sheep/goat|Schaf/Ziege",
  "game meat|Wild",
  "fish|Fisch",
  "undefined|unbestimmt"
)

$row = 218
foreach ($val in $species) {
  $ws.Cells.Item($row, 1).Value = "Species | Tierart "
  $ws.Cells.Item($row, 2).Value = $val
  $ws.Cells.Item($row, 2).Font.Name = "Arial"
  $ws.Cells.Item($row, 2).Font.Size = 10
  $ws.Cells.Item($row, 2).WrapText = $true
  $ws.Rows.Item($row).RowHeight = 15.05
  $row++
}

# --- column widths (small re-flow side effect of the new font / content) ---------------
$ws.Columns.Item(1).ColumnWidth = 63.8265306122449
$ws.Columns.Item(2).ColumnWidth = 73.8673469387755
$ws.Columns.Item(3).ColumnWidth = 35.530612244898

# --- selection / scroll position ---------------------------------------------------------
$ws.Range("B226").Select()
$excel.ActiveWindow.ScrollRow = 196
